# Update Desobhan (desbahan) yearly income-statement database:
#  - drop the oldest reported period (1396/12) and shift the remaining
#    periods one column to the left
#  - add a brand-new period column (1401/12) for which figures have not
#    been published yet
#  - refresh the "publish date" row with the newer disclosure dates
#  - the newly-shifted columns currently carry no confirmed figures yet
#    (pending the new read_price algorithm), so every financial figure is
#    reset to 0 (or "-" where the source row already used a placeholder)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 8: financial-period headers (D:H) ----
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# ---- Row 9: publish dates (D:H) ----
# H9 is a bare "yyyy-mm-dd" string (no trailing annotation), so it must be
# forced to text - otherwise Excel's COM layer auto-coerces it into a date
# serial number instead of keeping the literal label.
$ws.Range("D9").Value = "1399-03-21 (8)"
$ws.Range("E9").Value = "1400-03-02 (8)"
$ws.Range("F9").Value = "1401-03-08 (8)"
$ws.Range("G9").Value = "1402-02-28 (7)"
$ws.Range("H9").Value = "'1402-02-28"

# ---- Data rows: every reported figure resets to 0 for the shifted/new
#      columns, since the restated database has not been repopulated yet ----
$dataRows = @(11, 12, 13, 14, 16, 17, 18, 19, 20, 21, 22, 24, 25, 26, 27)
foreach ($r in $dataRows) {
    $ws.Range("D$r").Value = 0
    $ws.Range("E$r").Value = 0
    $ws.Range("F$r").Value = 0
    $ws.Range("G$r").Value = 0
    $ws.Range("H$r").Value = 0
}

# ---- Row 15 keeps its "-" placeholder in D, but E:G now become "-" too
#      (only H15 stays a reported 0) ----
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "-"
$ws.Range("H15").Value = 0

# ---- Row 23 becomes "-" across every period column ----
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = "-"
$ws.Range("H23").Value = "-"
